$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the two shapes we need to edit by name (more robust than a bare index).
$shp1 = $s.Shapes.Item("Rectangle 26")
$shp2 = $s.Shapes.Item("Rectangle 29")

# --- Edit 1: split "Query metrics from New Relic and Prometheus"
#     into three separate runs: "Query " / "metrics from New Relic and " / "Prometheus"
$tr1 = $shp1.TextFrame.TextRange
$tr1.Text = ""
$run1 = $tr1.InsertAfter("Query ")
$run2 = $tr1.InsertAfter("metrics from New Relic and ")
$run3 = $tr1.InsertAfter("Prometheus")

# --- Edit 2: collapse the triple space in
#     "Iter8 experiment   with A/B testing and progressive deployment"
$tr2 = $shp2.TextFrame.TextRange
$run = $tr2.Runs(2)
$run.Text = "Iter8 experiment with A/B testing and progressive deployment"
